$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row above row 41 (the old "TOTAL" row), shifting the TOTAL
# row and everything below it down by one. This gives us a new, blank data
# row 41 (project #37) while rows 42-46 keep the formatting that used to
# belong to rows 41-45.
# ---------------------------------------------------------------------------
$ws.Rows("41:41").Insert(-4121)

# ---------------------------------------------------------------------------
# New row 41 becomes the new last data row, taking on the thick-bottom-border
# formatting that row 40 (the current last data row) has right now.
# ---------------------------------------------------------------------------
$ws.Range("B40:K40").Copy()
$ws.Range("B41:K41").PasteSpecial(-4122)
$ws.Range("A41").Value = 37

# ---------------------------------------------------------------------------
# Row 40 used to be the last data row (thick bottom border). Now that a new
# data row follows it, it goes back to the regular "middle of table" look,
# matching the style used by rows 5-39.
# ---------------------------------------------------------------------------
$ws.Range("B39:K39").Copy()
$ws.Range("B40:K40").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# The TOTAL row (now row 42) needs its SUM/ratio formulas extended to cover
# the new row 41.
# ---------------------------------------------------------------------------
$ws.Range("D42").Formula = "=SUM(D5:D41)"
$ws.Range("E42").Formula = "=SUM(E5:E41)"
$ws.Range("F42").Formula = "=SUM(F5:F41)"
$ws.Range("G42").Formula = "=SUM(G5:G41)"
$ws.Range("H42").Formula = "=SUM(H5:H41)"
$ws.Range("I42").Formula = "=F42/H42"
$ws.Range("J42").Formula = "=F42/G42"
$ws.Range("K42").Formula = "=G42/H42"

# ---------------------------------------------------------------------------
# Update the print area to include the newly inserted row.
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintArea = "$B$1:$K$45"

# ---------------------------------------------------------------------------
# Match the saved selection state.
# ---------------------------------------------------------------------------
$ws.Range("B5").Select()
